# "last changes to v1.8.2"
# - Bump the IG version and publication date on the Metadata sheet.
# - Populate the missing ele-1/ext-1 invariant text on the "Extension" row
#   (AJ1) of the Elements sheet (it was previously only present on the
#   Element.extension row).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.2"
$meta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
